$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-20 20:18:44'
$ws.Range('N2').Value = '-2.3 °C 19:55 TU'
$ws.Range('E3').Value = '2026-02-20 20:18:46'
$ws.Range('E4').Value = '2026-02-20 20:18:49'
$ws.Range('H4').Value = '57%'
$ws.Range('J4').Value = '1022.5 hPa'
$ws.Range('K4').Value = '7.8 MJ/m2'
$ws.Range('O4').Value = '10.3 °C'
$ws.Range('E5').Value = '2026-02-20 20:18:51'
$ws.Range('E6').Value = '2026-02-20 20:18:54'
$ws.Range('J6').Value = '1022.5 hPa'
$ws.Range('O6').Value = '9.5 °C'
$ws.Range('E7').Value = '2026-02-20 20:18:56'
$ws.Range('H7').Value = '47%'
$ws.Range('J7').Value = '1022.4 hPa'
$ws.Range('E8').Value = '2026-02-20 20:18:59'
$ws.Range('H8').Value = '60%'
$ws.Range('J8').Value = '1022.7 hPa'
$ws.Range('O8').Value = '9.2 °C'
$ws.Range('E9').Value = '2026-02-20 20:19:02'
$ws.Range('E10').Value = '2026-02-20 20:19:05'
$ws.Range('O10').Value = '7.8 °C'
$ws.Range('E11').Value = '2026-02-20 20:19:07'
$ws.Range('H11').Value = '32%'
$ws.Range('O11').Value = '9.4 °C'
$ws.Range('E12').Value = '2026-02-20 20:19:10'
$ws.Range('E13').Value = '2026-02-20 20:19:13'
$ws.Range('J13').Value = '1023.5 hPa'
$ws.Range('O13').Value = '6.5 °C'
$ws.Range('E14').Value = '2026-02-20 20:19:16'
$ws.Range('H14').Value = '56%'
$ws.Range('O14').Value = '12.1 °C'
$ws.Range('E15').Value = '2026-02-20 20:19:18'
$ws.Range('E16').Value = '2026-02-20 20:19:21'
$ws.Range('O16').Value = '-3.2 °C'
$ws.Range('E17').Value = '2026-02-20 20:19:23'
$ws.Range('O17').Value = '3.0 °C'
$ws.Range('E18').Value = '2026-02-20 20:19:26'
$ws.Range('J18').Value = '1022.8 hPa'
$ws.Range('O18').Value = '8.0 °C'
$ws.Range('E19').Value = '2026-02-20 20:19:29'
$ws.Range('E20').Value = '2026-02-20 20:19:32'
$ws.Range('E21').Value = '2026-02-20 20:19:35'
$ws.Range('J21').Value = '1022.5 hPa'
$ws.Range('O21').Value = '9.4 °C'
$ws.Range('E22').Value = '2026-02-20 20:19:37'
$ws.Range('O22').Value = '-4.1 °C'
$ws.Range('E23').Value = '2026-02-20 20:19:40'
$ws.Range('H23').Value = '66%'
$ws.Range('O23').Value = '-4.9 °C'
$ws.Range('E24').Value = '2026-02-20 20:19:43'
$ws.Range('J24').Value = '1025.3 hPa'
$ws.Range('E25').Value = '2026-02-20 20:19:46'
$ws.Range('E26').Value = '2026-02-20 20:19:48'
$ws.Range('J26').Value = '1021.6 hPa'
$ws.Range('E27').Value = '2026-02-20 20:19:51'
$ws.Range('E28').Value = '2026-02-20 20:19:54'
$ws.Range('H28').Value = '64%'
$ws.Range('J28').Value = '1022.9 hPa'
$ws.Range('E29').Value = '2026-02-20 20:19:57'
$ws.Range('H29').Value = '70%'
$ws.Range('O29').Value = '9.8 °C'
$ws.Range('E30').Value = '2026-02-20 20:19:59'
$ws.Range('H30').Value = '58%'
$ws.Range('J30').Value = '1022.2 hPa'
$ws.Range('O30').Value = '11.1 °C'
$ws.Range('E31').Value = '2026-02-20 20:20:02'
$ws.Range('J31').Value = '1021.4 hPa'
$ws.Range('E32').Value = '2026-02-20 20:20:05'
$ws.Range('E33').Value = '2026-02-20 20:20:07'
$ws.Range('H33').Value = '41%'
$ws.Range('J33').Value = '1022.9 hPa'
$ws.Range('N33').Value = '2.6 °C 19:51 TU'
$ws.Range('O33').Value = '6.1 °C'
$ws.Range('E34').Value = '2026-02-20 20:20:10'
$ws.Range('H34').Value = '44%'
$ws.Range('L34').Value = '95.0 km/h - 20º 19:32 TU'
$ws.Range('O34').Value = '0.8 °C'
$ws.Range('E35').Value = '2026-02-20 20:20:13'
$ws.Range('J35').Value = '1026.7 hPa'
$ws.Range('E36').Value = '2026-02-20 20:20:16'
$ws.Range('J36').Value = '1022.4 hPa'
$ws.Range('E37').Value = '2026-02-20 20:20:18'
$ws.Range('H37').Value = '66%'
$ws.Range('J37').Value = '1024.3 hPa'
$ws.Range('E38').Value = '2026-02-20 20:20:21'
$ws.Range('O38').Value = '8.9 °C'
$ws.Range('E39').Value = '2026-02-20 20:20:24'
$ws.Range('O39').Value = '-2.6 °C'
$ws.Range('E40').Value = '2026-02-20 20:20:26'
$ws.Range('J40').Value = '1023.3 hPa'
$ws.Range('O40').Value = '10.5 °C'
$ws.Range('E41').Value = '2026-02-20 20:20:29'
$ws.Range('J41').Value = '1023.0 hPa'
$ws.Range('E42').Value = '2026-02-20 20:20:32'
$ws.Range('O42').Value = '10.3 °C'
$ws.Range('E43').Value = '2026-02-20 20:20:35'
$ws.Range('E44').Value = '2026-02-20 20:20:38'
$ws.Range('H44').Value = '77%'
$ws.Range('O44').Value = '-4.7 °C'
$ws.Range('E45').Value = '2026-02-20 20:20:40'
$ws.Range('J45').Value = '1029.5 hPa'
$ws.Range('O45').Value = '3.6 °C'
$ws.Range('E46').Value = '2026-02-20 20:20:43'
